{"js": "// The document has an item in the list that reads \"Fill out results\" but is\n// split across two runs (\"Fill out \" and \"results\") with identical\n// formatting. Normalize it into a single run containing \"Fill out results\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text === \"Fill out results\");\nif (!target) {\n  throw new Error('Could not find paragraph with text \"Fill out results\"');\n}\n\n// Re-insert the same visible text as a single run; Word coalesces it into\n// one run (since the paragraph mark / formatting carries over), replacing\n// the previous \"Fill out \" + \"results\" run pair.\ntarget.insertText(\"Fill out results\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The numbered list item that reads \"Fill out results\" is stored as two\n# adjacent runs (\"Fill out \" and \"results\") that share identical\n# character formatting. Normalize it to a single run with the text\n# \"Fill out results\".\n\n$d = $word.ActiveDocument\n\n# Find the specific paragraph (avoids touching \"Fill out study goals\" /\n# \"Fill out used methods\", which are separate list items).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Fill out results\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph with text 'Fill out results'\"\n}\n\n$r = $target.Range\n$r.MoveEnd(1, -1) | Out-Null  # exclude the paragraph mark from the range\n\n$find = $r.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Fill out results\"\n$find.Replacement.Text = \"Fill out results\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n"}
